$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (font, border, alignment) from the last existing
# header cell (G1) onto the new header cell (H1), then set its text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# Populate the new "Save" column data values
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 0
